# The "Project" and "Role" columns are being dropped from the daily-log
# template (storage/daily-activity-service import). On the original sheet
# the header row reads: B=Date, C=Project, D=What you've done today ?,
# E=Role, F=Duration, G=Any issue or note you want to share ?
#
# Removing the Project column (C) shifts everything left by one, so the
# old Role column is now D; deleting that brings Duration/Notes up into
# D/E and leaves the sheet as: B=Date, C=What you've done today ?,
# D=Duration, E=Any issue or note you want to share ?

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("Project"): remaining columns shift left.
$ws.Columns("C").Delete() | Out-Null

# After the shift, the former "Role" column (was E) is now D - delete it too.
$ws.Columns("D").Delete() | Out-Null

# Leave the selection where the edit ended up.
$ws.Range("C12").Select() | Out-Null
